# Tripadvisor New Orleans shard 195 update:
#   1. hotel_info gains a new "State" column (value "Louisiana") right
#      after "Hotel_Name" and before "City".
#   2. The two worksheets swap tab order: review_info moves to the front,
#      hotel_info moves to the back.

$wb = $excel.ActiveWorkbook

$wsHotel  = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# --- 1. Insert the new "State" column into hotel_info -----------------
# Do this while $wsHotel still unambiguously refers to the hotel_info
# worksheet (i.e. before any tab reordering below).
$wsHotel.Columns("C").Insert()
$wsHotel.Range("C1").Value = "State"
$wsHotel.Range("C2").Value = "Louisiana"

# --- 2. Reorder the tabs: review_info first, hotel_info second --------
$wsHotel.Move($null, $wsReview)
